$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell used as a style donor to restore formatting on cells that need a temporary
# Text number-format applied so Excel does not auto-convert numeric-looking strings
# (e.g. "7.20", "0.999") into floating point numbers, which would silently drop
# significant trailing zeros / change the stored type.
$plainStyle = $ws.Range("D4").Style

$ws.Range("D2").Value = "64.332.57"
$ws.Range("E2").Value = "  +0.58%  "

$ws.Range("D3").Value = "3.495.89"
$ws.Range("E3").Value = "  -0.30%  "

$ws.Range("E4").Value = "  -0.16%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "587.41"
$ws.Range("D5").Style = $plainStyle
$ws.Range("E5").Value = "  +0.35%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "134.44"
$ws.Range("D6").Style = $plainStyle
$ws.Range("E6").Value = "  +1.69%  "

$ws.Range("D7").Value = "3.495.90"
$ws.Range("E7").Value = "  -0.27%  "

$ws.Range("E8").Value = "  -0.04%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.486"
$ws.Range("D9").Style = $plainStyle

$ws.Range("E10").Value = "  -0.34%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "7.20"
$ws.Range("D11").Style = $plainStyle
$ws.Range("E11").Value = "  +1.29%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.377"
$ws.Range("D12").Style = $plainStyle
$ws.Range("E12").Value = "  -2.31%  "

$ws.Range("D13").Value = "4.084.27"
$ws.Range("E13").Value = "  -0.70%  "

$ws.Range("E14").Value = "  +1.82%  "

$ws.Range("E15").Value = "  +0.25%  "

$ws.Range("D16").Value = "3.493.04"

$ws.Range("D17").Value = "64.358.50"
$ws.Range("E17").Value = "  +0.41%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "25.21"
$ws.Range("D18").Style = $plainStyle
$ws.Range("E18").Value = "  -9.33%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "10.02"
$ws.Range("D19").Style = $plainStyle
$ws.Range("E19").Value = "  +0.24%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "5.72"
$ws.Range("D20").Style = $plainStyle
$ws.Range("E20").Value = "  +1.00%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.65"
$ws.Range("D21").Style = $plainStyle
$ws.Range("E21").Value = "  -5.46%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "388.34"
$ws.Range("D22").Style = $plainStyle
$ws.Range("E22").Value = "  -0.64%  "

$ws.Range("E23").Value = "  -2.13%  "

$ws.Range("D24").Value = "3.634.02"
$ws.Range("E24").Value = "  -0.51%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "74.25"
$ws.Range("D25").Style = $plainStyle
$ws.Range("E25").Value = "  +1.77%  "

$ws.Range("E26").Value = "  +0.04%  "

$ws.Range("E27").Value = "  -0.21%  "

$ws.Range("E28").Value = "  +1.02%  "

$ws.Range("B29").Value = "Fetch.AI"
$ws.Range("C29").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.54"
$ws.Range("D29").Style = $plainStyle
$ws.Range("E29").Value = "  -2.18%  "

$ws.Range("B30").Value = "Binance-PegBSC-USD"
$ws.Range("C30").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.999"
$ws.Range("D30").Style = $plainStyle
$ws.Range("E30").Value = "  -0.07%  "

$ws.Range("E31").Value = "  -0.59%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "8.30"
$ws.Range("D32").Style = $plainStyle
$ws.Range("E32").Value = "  +0.39%  "

$ws.Range("E33").Value = "  -1.24%  "

$ws.Range("D34").Value = "3.514.46"
$ws.Range("E34").Value = "  +0.11%  "

$ws.Range("E36").Value = "  +2.40%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "23.44"
$ws.Range("D37").Style = $plainStyle
$ws.Range("E37").Value = "  -2.26%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "5.26"
$ws.Range("D38").Style = $plainStyle
$ws.Range("E38").Value = "  -1.91%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "6.85"
$ws.Range("D39").Style = $plainStyle
$ws.Range("E39").Value = "  -2.16%  "

$ws.Range("E40").Value = "  -2.17%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "161.16"
$ws.Range("D41").Style = $plainStyle
$ws.Range("E41").Value = "  -4.24%  "

$ws.Range("E42").Value = "  -3.83%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.806"
$ws.Range("D43").Style = $plainStyle
$ws.Range("E43").Value = "  -0.87%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "25.57"
$ws.Range("D44").Style = $plainStyle
$ws.Range("E44").Value = "  -4.25%  "

$ws.Range("E45").Value = "  -0.19%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "41.80"
$ws.Range("D46").Style = $plainStyle
$ws.Range("E46").Value = "  -0.14%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "4.41"
$ws.Range("D47").Style = $plainStyle
$ws.Range("E47").Value = "  +0.59%  "

$ws.Range("E48").Value = "  -0.70%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.65"
$ws.Range("D49").Style = $plainStyle
$ws.Range("E49").Value = "  +0.59%  "

$ws.Range("D50").Value = "2.472.62"
$ws.Range("E50").Value = "  +0.93%  "

$ws.Range("E51").Value = "  -2.13%  "
